# Apply the task_summary.xlsx refresh: bump the check/start/end timestamps
# and the task-date-derived names forward a few days (re-run of the same
# report), and fill in a previously-missing Start_time for the HSG task
# (row 9) that had errored out before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (UZH LARGE READY) - refreshed check time only.
$ws.Range("C2").Value = "2024-07-24 11:41:01"

# Row 8 (ZBZ SMALL READY) - task re-scheduled from 2024-07-22 to 2024-07-26.
$ws.Range("B8").Value = "task_2024-07-26_ZBZ_SMALL_READY"
$ws.Range("C8").Value = "2024-07-24 11:41:55"
# Leading apostrophe forces text so the date-shaped string isn't reinterpreted
# as a serial date number (column holds plain text dates, e.g. "2024-07-26").
$ws.Range("F8").Value = "'2024-07-26"

# Row 9 (HSG SMALL DONE) - task re-run from 2024-07-20 to 2024-07-24, and
# this run recorded a Start_time (previously blank, i.e. the task had
# errored before starting).
$ws.Range("B9").Value = "task_2024-07-24_HSG_SMALL_DONE"
$ws.Range("C9").Value = "2024-07-24 11:42:05"
$ws.Range("D9").Value = "2024-07-24 11:42:08"
$ws.Range("E9").Value = "2024-07-24 11:42:15"
$ws.Range("F9").Value = "'2024-07-24"
